$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46075 -> 46076, i.e. 2026-02-22 -> 2026-02-23) for every data row (2..332).
for ($r = 2; $r -le 332; $r++) {
    $ws.Cells.Item($r, 3).Value = 46076
}
